$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44965
$ws.Range("J2").Value = 87
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 3000
$ws.Range("P2").Value = 1000

# Row 3
$ws.Range("D3").Value = 44223
$ws.Range("J3").Value = 80
$ws.Range("M3").Value = 2781
$ws.Range("P3").Value = 927

# Row 4
$ws.Range("D4").Value = 44225
$ws.Range("J4").Value = 56
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("P4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44179
$ws.Range("J5").Value = 78
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("P5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44559
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 667

# Row 7
$ws.Range("D7").Value = 44574
$ws.Range("J7").Value = 50

# Row 8
$ws.Range("D8").Value = 45189
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44804
$ws.Range("J9").Value = 85

# Row 10
$ws.Range("D10").Value = 44967
$ws.Range("J10").Value = 110
$ws.Range("L10").Value = 3300
$ws.Range("M10").Value = 3136
$ws.Range("P10").Value = 1045

# Row 11
$ws.Range("D11").Value = 45092
$ws.Range("J11").Value = 90
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = 3278
$ws.Range("P11").Value = 1093

# Row 12
$ws.Range("D12").Value = 44669
$ws.Range("J12").Value = 92
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 2755
$ws.Range("P12").Value = 918

# Row 13
$ws.Range("D13").Value = 44627
$ws.Range("J13").Value = 78
$ws.Range("K13").Value = 3500
$ws.Range("L13").Value = 3500
$ws.Range("M13").Value = 3500
$ws.Range("P13").Value = 1167

# Row 14
$ws.Range("D14").Value = 44291
$ws.Range("J14").Value = 45

# Row 15
$ws.Range("D15").Value = 44756
$ws.Range("J15").Value = 104
$ws.Range("K15").Value = 2800
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 2904
$ws.Range("P15").Value = 968

# Row 16
$ws.Range("D16").Value = 44390
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("P16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44845
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2500
$ws.Range("P17").Value = 833

# Row 18
$ws.Range("D18").Value = 44221
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2500
$ws.Range("P18").Value = 833

# Row 19
$ws.Range("D19").Value = 44224
$ws.Range("J19").Value = 67

# Row 20
$ws.Range("D20").Value = 44260
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 3500
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = 3500
$ws.Range("P20").Value = 1167

# Row 21
$ws.Range("D21").Value = 44166
$ws.Range("J21").Value = 45
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 833

# Row 22
$ws.Range("D22").Value = 44193
$ws.Range("J22").Value = 70

# Row 23
$ws.Range("D23").Value = 44992
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 4000
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = 4000
$ws.Range("P23").Value = 1333

# Row 24
$ws.Range("D24").Value = 44537
$ws.Range("J24").Value = 88
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2200
$ws.Range("M24").Value = 2091
$ws.Range("P24").Value = 697

# Row 25
$ws.Range("D25").Value = 44243

# Row 26
$ws.Range("D26").Value = 44340
$ws.Range("J26").Value = 54
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = 3000
$ws.Range("P26").Value = 1000

# Row 27
$ws.Range("D27").Value = 45118
$ws.Range("J27").Value = 67
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = 3000
$ws.Range("P27").Value = 1000

# Row 28
$ws.Range("D28").Value = 44389
$ws.Range("J28").Value = 81
$ws.Range("K28").Value = 2800
$ws.Range("M28").Value = 2889
$ws.Range("P28").Value = 963

# Row 29
$ws.Range("D29").Value = 44935
$ws.Range("J29").Value = 78
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 3000
$ws.Range("P29").Value = 1000

# Row 30
$ws.Range("D30").Value = 44536
$ws.Range("J30").Value = 125
$ws.Range("K30").Value = 2200
$ws.Range("L30").Value = 2200
$ws.Range("M30").Value = 2200
$ws.Range("P30").Value = 733

# Row 31
$ws.Range("D31").Value = 44937
$ws.Range("J31").Value = 68
$ws.Range("K31").Value = 3500
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = 3500
$ws.Range("P31").Value = 1167

# Row 32
$ws.Range("D32").Value = 44187
$ws.Range("J32").Value = 65
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = 3000
$ws.Range("P32").Value = 1000

# Row 33
$ws.Range("D33").Value = 45117
$ws.Range("J33").Value = 56

# Row 34
$ws.Range("D34").Value = 44292
$ws.Range("J34").Value = 40
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = 3000
$ws.Range("P34").Value = 1000

# Row 35
$ws.Range("D35").Value = 44222
$ws.Range("J35").Value = 45
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = 3000
$ws.Range("P35").Value = 1000

# Row 36
$ws.Range("D36").Value = 44165
$ws.Range("J36").Value = 68
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = 3000
$ws.Range("P36").Value = 1000

# Row 37
$ws.Range("D37").Value = 44242
$ws.Range("J37").Value = 95
$ws.Range("K37").Value = 2500
$ws.Range("M37").Value = 2737
$ws.Range("P37").Value = 912

# Row 38
$ws.Range("D38").Value = 44557
$ws.Range("J38").Value = 104
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 2500
$ws.Range("M38").Value = 2260
$ws.Range("P38").Value = 753
